$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = "Мы все еще находимся в процессе приобретения товара, но вы можете сделать предварительный заказ, чтобы мы знали, что вы заинтересованы. Мы свяжемся с вами, когда товар будет в наличии, чтобы договориться об оплате."
$ws.Range("A25").Value = "We are still in the process of acquiring stock, but you are welcome to pre-order so we know that you are interested. We are going to contact you when stock is available to arrange payment."
$ws.Range("C25").Value = "We are still in the process of acquiring stock, but you are welcome to pre-order so we know that you are interested. We are going to contact you when stock is available to arrange payment."

$ws.Range("C25").Select()
